# Daily attendance processing - reorder "Recorded By" (column G) entries so
# that email-address entries come first, followed by non-email entries
# (e.g. "System" / "system"), preserving the relative order within each group.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $val = $cell.Value2

    if ($null -eq $val -or $val -eq "") {
        continue
    }

    $parts = $val -split ','
    $trimmed = @()
    foreach ($p in $parts) {
        $trimmed += $p.Trim()
    }

    $emails = @()
    $others = @()
    foreach ($p in $trimmed) {
        if ($p -like '*@*') {
            $emails += $p
        } else {
            $others += $p
        }
    }

    $newOrder = $emails + $others
    $newVal = [string]::Join(', ', $newOrder)

    if ($newVal -ne $val) {
        $cell.Value = $newVal
    }
}
